function Set-MatchRow($Row, $Indice, $Data, $Home, $HomeGols, $Away, $AwayGols, $HomeOpenOdds, $HomeOpenHora, $HomeCloseOdds, $HomeCloseHora, $DrawOpenOdds, $DrawOpenHora, $DrawCloseOdds, $DrawCloseHora, $AwayOpenOdds, $AwayOpenHora, $AwayCloseOdds, $AwayCloseHora, $Url) {
    $ws.Cells.Item($Row, 1).Value = $Indice
    $ws.Cells.Item($Row, 2).Value = "moldova"
    $ws.Cells.Item($Row, 3).Value = "super-liga"
    $ws.Cells.Item($Row, 4).Value = "2023-2024"
    $ws.Cells.Item($Row, 5).Value2 = $Data
    $ws.Cells.Item($Row, 6).Value = $Home
    $ws.Cells.Item($Row, 7).Value = $HomeGols
    $ws.Cells.Item($Row, 8).Value = $Away
    $ws.Cells.Item($Row, 9).Value = $AwayGols
    $ws.Cells.Item($Row, 10).Value = $HomeOpenOdds
    $ws.Cells.Item($Row, 11).Value = $HomeOpenHora
    $ws.Cells.Item($Row, 12).Value = $HomeCloseOdds
    $ws.Cells.Item($Row, 13).Value = $HomeCloseHora
    $ws.Cells.Item($Row, 14).Value = $DrawOpenOdds
    $ws.Cells.Item($Row, 15).Value = $DrawOpenHora
    $ws.Cells.Item($Row, 16).Value = $DrawCloseOdds
    $ws.Cells.Item($Row, 17).Value = $DrawCloseHora
    $ws.Cells.Item($Row, 18).Value = $AwayOpenOdds
    $ws.Cells.Item($Row, 19).Value = $AwayOpenHora
    $ws.Cells.Item($Row, 20).Value = $AwayCloseOdds
    $ws.Cells.Item($Row, 21).Value = $AwayCloseHora
    $ws.Cells.Item($Row, 22).Value = $Url
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 40 and 41 were reordered (the match that used to be listed second
#    is now listed first). Swap the match-data columns F:V between the two
#    rows; columns A (row index) and E (match date) stay put.
# ---------------------------------------------------------------------------
$rngA = $ws.Range("F40:V40")
$rngB = $ws.Range("F41:V41")
$valsA = $rngA.Value2
$valsB = $rngB.Value2
$rngA.Value2 = $valsB
$rngB.Value2 = $valsA

# ---------------------------------------------------------------------------
# 2) Three new match rows (50-52) were appended at the end of the sheet.
#    Copy the formatting of the last existing data row (49) down into the
#    new rows first, so the new cells pick up the same styles (bordered /
#    centered index column, date-time number format for the match date
#    column, etc.), then overwrite the copied values with the real data.
# ---------------------------------------------------------------------------
$lastRow = $ws.Range("A49:V49")
$lastRow.Copy($ws.Range("A50:V50"))
$lastRow.Copy($ws.Range("A51:V51"))
$lastRow.Copy($ws.Range("A52:V52"))

Set-MatchRow 50 49 45255.5 "Milsami" 3 "Sparta Selemet" 1 1.17 "25/11/2023 01:13" 1.06 "25/11/2023 11:56" 6.28 "25/11/2023 01:13" 9.9 "25/11/2023 11:59" 10.24 "25/11/2023 01:13" 18.06 "25/11/2023 11:59" "https://www.betexplorer.com/football/moldova/super-liga/milsami-sparta-selemet/YqNtwrYl/"

Set-MatchRow 51 50 45255.54166666666 "Sheriff Tiraspol" 2 "Zimbru Chisinau" 0 1.28 "25/11/2023 02:12" 1.31 "25/11/2023 12:54" 4.3 "25/11/2023 02:12" 4.32 "25/11/2023 12:54" 10.37 "25/11/2023 02:12" 9.25 "25/11/2023 12:54" "https://www.betexplorer.com/football/moldova/super-liga/sheriff-tiraspol-zimbru-chisinau/jLXyvOIr/"

Set-MatchRow 52 51 45255.625 "Petrocub" 3 "Balti" 0 1.64 "25/11/2023 04:13" 1.44 "25/11/2023 11:54" 3.56 "25/11/2023 04:13" 3.98 "25/11/2023 14:47" 4.32 "25/11/2023 04:13" 5.42 "25/11/2023 11:54" "https://www.betexplorer.com/football/moldova/super-liga/petrocub-hincesti-csf-balti/CzMpx2me/"
